# "one giant update whoohooooo"
# Update the daily diary notes on the "Projects Overview" sheet and fix the
# weekend highlighting that had drifted onto the wrong day.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Projects Overview")

# --- Monday 12 June (row 31): now marked as "done" (green) and gets its first note ---
$ws.Range("C25").Copy() | Out-Null
$ws.Range("C31").PasteSpecial(-4122) | Out-Null
$ws.Range("D31").Value = "damian heeft poster's gemaakt en alex maakt de app"

# --- Tuesday 13 June (row 32): done (green) + note ---
$ws.Range("C25").Copy() | Out-Null
$ws.Range("C32").PasteSpecial(-4122) | Out-Null
$ws.Range("D32").Value = "damian heeft adjustments gemaakt aan de posters en alex werkt aan de app"

# --- Wednesday 14 June (row 33): done (green) + note ---
$ws.Range("C25").Copy() | Out-Null
$ws.Range("C33").PasteSpecial(-4122) | Out-Null
$ws.Range("D33").Value = "alex werkt og aan de app en damian maakt user inventory"

# --- Thursday 15 June (row 34): done (green) + same note repeated ---
$ws.Range("C25").Copy() | Out-Null
$ws.Range("C34").PasteSpecial(-4122) | Out-Null
$ws.Range("D34").Value = "alex werkt og aan de app en damian maakt user inventory"

# --- Saturday 17 June (row 36): this is the real weekend day, so flag it ---
$ws.Range("C29").Copy() | Out-Null
$ws.Range("C36").PasteSpecial(-4122) | Out-Null
$ws.Range("D29").Copy() | Out-Null
$ws.Range("D36").PasteSpecial(-4122) | Out-Null
$ws.Range("D36").Value = "weekend"
$ws.Range("H29").Copy() | Out-Null
$ws.Range("H36").PasteSpecial(-4122) | Out-Null

# --- Monday 19 June (row 38): was mislabeled as weekend, revert to a normal day ---
$ws.Range("C39").Copy() | Out-Null
$ws.Range("C38").PasteSpecial(-4122) | Out-Null
$ws.Range("D38").Clear() | Out-Null
$ws.Range("H38").Clear() | Out-Null

$excel.CutCopyMode = 0

# Restore the selection to where the author ended up after making the edits.
$ws.Range("D9").Select() | Out-Null
